# Daily attendance processing - 2026-01-22 08:44:19
#
# Applies the following updates to the "Session Analysis Results" sheet:
#   1. Swap the "Recorded By" text order for sessions recorded by both the
#      System and the user (column G) from "System, dnasr281@gmail.com" to
#      "dnasr281@gmail.com, System".
#   2. Refresh the dashboard summary counters in K7:L8
#      (Missing Sessions / Pending Sessions).
#   3. Refresh the per-group Pending/Missing session counters (columns P:Q)
#      for the affected groups (rows 21-26).
#   4. Re-classify six sessions (rows 184, 211, 238, 265, 292, 319) whose
#      status was stuck at "Pending" -> "Not Recorded", matching the
#      formatting/status convention used elsewhere in the sheet (style of
#      row 21, i.e. the pink "Not Recorded" highlight instead of the
#      yellow "Pending" highlight).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

# --- 1. Swap "System, <email>" -> "<email>, System" in column G --------------
$gRows = @(8, 9, 10, 12, 14, 15, 17, 34, 35, 36, 38, 40, 41, 43, 60, 61, 62, 64, 66, 67, 69, 86, 87, 88, 90, 92, 93, 95, 112, 113, 114, 116, 118, 119, 121, 138, 139, 140, 142, 144, 145, 147)
foreach ($r in $gRows) {
    $ws.Cells.Item($r, 7).Value = "dnasr281@gmail.com, System"
}

# --- 2. Dashboard summary counters -------------------------------------------
$ws.Range("L7").Value = 69
$ws.Range("L8").Value = 0

# --- 3. Per-group Pending / Missing counters (columns P & Q) -----------------
$ws.Range("P21").Value = 7
$ws.Range("Q21").Value = 0
$ws.Range("P22").Value = 7
$ws.Range("Q22").Value = 0
$ws.Range("P23").Value = 7
$ws.Range("Q23").Value = 0
$ws.Range("P24").Value = 8
$ws.Range("Q24").Value = 0
$ws.Range("P25").Value = 7
$ws.Range("Q25").Value = 0
$ws.Range("P26").Value = 7
$ws.Range("Q26").Value = 0

# --- 4. Re-classify "Pending" sessions as "Not Recorded" ---------------------
# Copy the formatting used by the existing "Not Recorded" row (21, columns
# A:I) onto each affected row, then update the Status text in column I.
$formatSourceRange = $ws.Range("A21:I21")
$formatSourceRange.Copy() | Out-Null

$pendingRows = @(184, 211, 238, 265, 292, 319)
foreach ($r in $pendingRows) {
    $targetRange = $ws.Range("A" + $r + ":I" + $r)
    $targetRange.PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    $ws.Cells.Item($r, 9).Value = "Not Recorded"
}

$excel.CutCopyMode = 0
